# Update cryptos list with latest scraped price/volume values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.252.09"
$ws.Range("E2").Value = "  +0.84%  "
$ws.Range("D3").Value = "1.822.26"
$ws.Range("E3").Value = "  +0.00%  "
$ws.Range("D4").Value = "'1.000"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").Value = "'313.19"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.69%  "
$ws.Range("D6").Value = "'1.001"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.08%  "
$ws.Range("D7").Value = "'0.4481"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.18%  "
$ws.Range("E8").Value = "  +1.72%  "
$ws.Range("D9").Value = "'0.07401"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.25%  "
$ws.Range("D10").Value = "'0.8789"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.50%  "
$ws.Range("D11").Value = "'20.84"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Value = "1.816.39"
$ws.Range("E12").Value = "  -0.15%  "
$ws.Range("D13").Value = "'6.713"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.16%  "
$ws.Range("D14").Value = "'5.415"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.56%  "
$ws.Range("D15").Value = "'92.82"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.64%  "
$ws.Range("D16").Value = "'0.07060"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.58%  "
$ws.Range("D17").Value = "'1.001"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").Value = "'0.000008803"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.24%  "
$ws.Range("D19").Value = "'1.001"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.01%  "
$ws.Range("D20").Value = "'15.02"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.36%  "
$ws.Range("D21").Value = "27.245.05"
$ws.Range("E21").Value = "  +0.70%  "
$ws.Range("D22").Value = "'5.341"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.36%  "
$ws.Range("D23").Value = "'10.94"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.19%  "
$ws.Range("E24").Value = "  -1.44%  "
$ws.Range("D25").Value = "'151.00"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.53%  "
$ws.Range("D26").Value = "'2.280"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.24%  "
$ws.Range("E27").Value = "  +0.42%  "
$ws.Range("D28").Value = "'5.341"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.63%  "
$ws.Range("D29").Value = "'117.34"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.57%  "
$ws.Range("D30").Value = "'0.08885"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.23%  "
$ws.Range("D31").Value = "'0.7890"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +4.49%  "
$ws.Range("D32").Value = "'1.196"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.08%  "
$ws.Range("D33").Value = "'4.572"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.66%  "
$ws.Range("D34").Value = "'2.926"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.21%  "
$ws.Range("E35").Value = "  -0.04%  "
$ws.Range("D36").Value = "'1.107"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.49%  "
$ws.Range("D37").Value = "'0.01974"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").Value = "'0.05259"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.50%  "
$ws.Range("D39").Value = "'7.286"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.27%  "
$ws.Range("D40").Value = "'0.5287"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.60%  "
$ws.Range("D41").Value = "'2.872"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.53%  "
$ws.Range("D42").Value = "'2.326"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +17.37%  "
$ws.Range("E43").Value = "  -0.31%  "
$ws.Range("D44").Value = "'8.630"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.49%  "
$ws.Range("D45").Value = "'0.5039"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.69%  "
$ws.Range("D46").Value = "'10.61"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.31%  "
$ws.Range("D47").Value = "'105.22"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.30%  "
$ws.Range("D48").Value = "'1.685"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").Value = "'1.000"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.04%  "
$ws.Range("D50").Value = "'0.06382"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.08%  "
$ws.Range("D51").Value = "'66.00"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +4.67%  "
